$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("best strategy compiler")

# Remove the old data (rows 3:52, cols A:E) so stale cells beyond the
# new table (row 33) don't linger.
$ws.Range("A3:E52").ClearContents()

# New header row describing the Strategy struct fields.
$ws.Range("A2").Value = "Burn"
$ws.Range("B2").Value = "BHoles"
$ws.Range("C2").Value = "FHoles"
$ws.Range("D2").Value = "HighY"
$ws.Range("E2").Value = "Step"

# New strategy parameter rows (Burn, BHoles, FHoles, HighY, Step).
$data = @(
  @(3,  3,7,3,2,1),
  @(4,  5,7,4,1,2),
  @(5,  1,6,4,2,1),
  @(6,  2,5,3,2,2),
  @(7,  4,6,4,3,1),
  @(8,  1,4,3,1,1),
  @(9,  2,4,3,1,1),
  @(10, 1,6,4,1,2),
  @(11, 3,6,3,1,2),
  @(12, 3,7,4,3,2),
  @(13, 5,4,2,1,1),
  @(14, 1,6,4,3,1),
  @(15, 2,6,4,1,1),
  @(16, 2,7,2,1,2),
  @(17, 3,7,4,1,2),
  @(18, 4,7,3,4,2),
  @(19, 1,5,2,2,1),
  @(20, 1,7,4,2,2),
  @(21, 2,7,4,2,2),
  @(22, 5,4,4,1,2),
  @(23, 4,6,3,1,1),
  @(24, 3,7,2,4,1),
  @(25, 3,7,2,4,2),
  @(26, 4,7,4,2,2),
  @(27, 5,7,3,2,2),
  @(28, 2,5,2,1,2),
  @(29, 2,5,4,2,1),
  @(30, 4,5,2,1,1),
  @(31, 1,7,4,3,2),
  @(32, 3,7,3,1,2),
  @(33, 4,7,2,3,1)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r,1).Value = $row[1]
  $ws.Cells.Item($r,2).Value = $row[2]
  $ws.Cells.Item($r,3).Value = $row[3]
  $ws.Cells.Item($r,4).Value = $row[4]
  $ws.Cells.Item($r,5).Value = $row[5]
}

# Column G now holds the generated "go playGames(...)" strategy call.
# Row 3 gets a standalone formula; rows 4:33 become a shared-formula
# group (mirrors the original layout's E3 standalone / E4:E52 shared).
$ws.Range("G3").Formula = '=CONCATENATE("go playGames(Strategy{Burn: ",A3,", Step:",E3,", BHoles:",B3,", FHoles:",C3,", HighY:",D3,"}, 22, false, false)")'
$ws.Range("G4:G33").Formula = '=CONCATENATE("go playGames(Strategy{Burn: ",A4,", Step:",E4,", BHoles:",B4,", FHoles:",C4,", HighY:",D4,"}, 22, false, false)")'

# Column layout: drop the old column E width and widen the new column G.
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(7).ColumnWidth = 80.8

# Selection / view now centers on the new formula column.
[void]$ws.Range("G3:G33").Select()

# Sheet2: the saved view no longer pins a frozen top-left cell.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Application.ActiveWindow.ScrollRow = 1
$ws2.Application.ActiveWindow.ScrollColumn = 1

# Keep "visual" as the active/selected tab, as it was before the edit.
$wsVisual = $wb.Worksheets.Item("visual")
[void]$wsVisual.Activate()
